$d = $word.ActiveDocument

# 1) Update the visible text of the two outer paragraphs.
$d.Content.Find.Execute("Outer para 1", $true, $false, $false, $false, $false,
                         $true, 1, $false, "outer, before sect break", 2)
$d.Content.Find.Execute("Outer para 2", $true, $false, $false, $false, $false,
                         $true, 1, $false, "outer, after sect break", 2)

# 2) Insert a continuous section break right after the first paragraph's
#    text, then delete the paragraph mark that used to end that first
#    paragraph. That merges "outer, before sect break" into the paragraph
#    that now owns the new sectPr, so the sectPr ends up inside that
#    paragraph's pPr -- exactly as a real section break would look.
$p1 = $d.Paragraphs(1)
$breakPos = $p1.Range.End - 1
$r = $d.Range($breakPos, $breakPos)
$r.InsertBreak(3) # wdSectionBreakContinuous
$mark = $d.Range($breakPos, $breakPos + 1)
$mark.Delete()

# 3) Configure the page setup for the newly created (first) section -- the
#    one whose sectPr now lives inside the first paragraph's pPr.
$sec1 = $d.Sections(1)
$sec1.PageSetup.PageWidth = 595.3
$sec1.PageSetup.PageHeight = 841.9
$sec1.PageSetup.TopMargin = 70.85
$sec1.PageSetup.RightMargin = 70.85
$sec1.PageSetup.BottomMargin = 70.85
$sec1.PageSetup.LeftMargin = 70.85
$sec1.PageSetup.HeaderDistance = 35.4
$sec1.PageSetup.FooterDistance = 35.4
$sec1.PageSetup.Gutter = 0
$sec1.PageSetup.TextColumns.Spacing = 35.4

$hdr1 = $sec1.Headers(1)
$hdr1.PageNumbers.RestartNumberingAtSection = $true
$hdr1.PageNumbers.StartingNumber = 1

# 4) Configure the page setup for the outer (last) section, and mark it as
#    a continuous section (matches the body-level sectPr).
$sec2 = $d.Sections(2)
$sec2.PageSetup.SectionStart = 0 # wdSectionContinuous
$sec2.PageSetup.PageWidth = 595.3
$sec2.PageSetup.PageHeight = 841.9
$sec2.PageSetup.TopMargin = 70.85
$sec2.PageSetup.RightMargin = 70.85
$sec2.PageSetup.BottomMargin = 70.85
$sec2.PageSetup.LeftMargin = 70.85
$sec2.PageSetup.HeaderDistance = 35.4
$sec2.PageSetup.FooterDistance = 35.4
$sec2.PageSetup.Gutter = 0
$sec2.PageSetup.TextColumns.Spacing = 35.4
